# The deck currently carries the "Integral" theme on the slide master
# (ppt/theme/theme1.xml) and a plain "Office Theme" on the notes master
# (ppt/theme/theme2.xml). The target edit swaps these two themes: the
# slide master should end up with the plain Office colours, and the
# notes master should end up with the Integral colours.
#
# Font scheme and format scheme (fills/lines/effects) are identical
# between the two themes, only the 12 theme colours (and the cosmetic
# "name" attributes, which aren't exposed through automation) differ,
# so the swap is performed by re-writing each ThemeColorScheme entry.

$p = $ppt.ActivePresentation

# VBA-style RGB() packs a colour as r + g*256 + b*65536.
function VbaRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Current slide-master ("Integral") theme colours, in
# ThemeColorScheme.Item(1..12) order: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.
$integralColors = @(
    (VbaRGB 0x00 0x00 0x00),  # dk1
    (VbaRGB 0xFF 0xFF 0xFF),  # lt1
    (VbaRGB 0x45 0x5F 0x51),  # dk2
    (VbaRGB 0xE3 0xDE 0xD1),  # lt2
    (VbaRGB 0x99 0xCB 0x38),  # accent1
    (VbaRGB 0x63 0xA5 0x37),  # accent2
    (VbaRGB 0xE6 0xD0 0x24),  # accent3
    (VbaRGB 0xCC 0x97 0x00),  # accent4
    (VbaRGB 0x4E 0xB3 0xCF),  # accent5
    (VbaRGB 0x37 0x8D 0xA6),  # accent6
    (VbaRGB 0x6B 0x9F 0x25),  # hlink
    (VbaRGB 0xB2 0x6B 0x02)   # folHlink
)

# Current notes-master ("Office Theme") theme colours, same order.
$officeColors = @(
    (VbaRGB 0x00 0x00 0x00),  # dk1
    (VbaRGB 0xFF 0xFF 0xFF),  # lt1
    (VbaRGB 0x44 0x54 0x6A),  # dk2
    (VbaRGB 0xE7 0xE6 0xE6),  # lt2
    (VbaRGB 0x5B 0x9B 0xD5),  # accent1
    (VbaRGB 0xED 0x7D 0x31),  # accent2
    (VbaRGB 0xA5 0xA5 0xA5),  # accent3
    (VbaRGB 0xFF 0xC0 0x00),  # accent4
    (VbaRGB 0x44 0x72 0xC4),  # accent5
    (VbaRGB 0x70 0xAD 0x47),  # accent6
    (VbaRGB 0x05 0x63 0xC1),  # hlink
    (VbaRGB 0x95 0x4F 0x72)   # folHlink
)

# Slide master picks up the Office colours (theme1.xml -> Office Theme).
$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}

# Notes master picks up the Integral colours (theme2.xml -> Integral).
$notesScheme = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $integralColors[$i - 1]
}
